$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("English and Communication", 5417, "Female", 2001, "Intermediate_2"),
    @("French", 701, "Female", 2001, "Intermediate_2"),
    @("Mathematics", 5886, "Female", 2001, "Intermediate_2"),
    @("Biology", 2418, "Female", 2001, "Intermediate_2"),
    @("Chemistry", 775, "Female", 2001, "Intermediate_2"),
    @("Physics", 454, "Female", 2001, "Intermediate_2"),
    @("Computing", 522, "Female", 2001, "Intermediate_2")
)

$row = 9
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}
